$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.073.10'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.783.92'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +7.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '426.36'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +9.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.27'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +12.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.611'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +4.96%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.739'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +9.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.158'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000327'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.13'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +12.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.65'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +17.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.390.16'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +7.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.04'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +18.94%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.812.37'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +8.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.15'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +8.53%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +12.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '66.191.31'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '411.73'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.18'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +9.88%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +15.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.71'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '37.13'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +11.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.86'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +47.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.29'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +11.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.83'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +14.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.41'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.99'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +19.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '708.27'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +6.87%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +17.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.71'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '39.90'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +9.90%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.77'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +41.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.152'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.03'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0474'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +9.63%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +51.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0687'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +13.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.88'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +8.85%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.52%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +8.01%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +11.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.321'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +18.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.17'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.90%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +6.83%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.64'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.99%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '143.13'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.81'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.20%  '
